# Exposures_DownloadAccountDetails.xlsx - documentation/notes update
# "More documentation and notes, No API changes."
#
# Summary of content edits:
#  - About sheet: "Last Update" value May 2023 -> May 2024
#  - REQUEST sheet (Table1): rename "Account Type"/"Account Number"/
#    "Account Check Digit" to the new Entity-* wording, and add new
#    "Model" notes (Model/optional) for two of the rows
#  - RESPONSE sheet (Table2): same renames as REQUEST for the Description
#    column entries

$wb = $excel.ActiveWorkbook

$wsAbout    = $wb.Worksheets.Item("About")
$wsRequest  = $wb.Worksheets.Item("REQUEST")
$wsResponse = $wb.Worksheets.Item("RESPONSE")

# --- About sheet: update "Last Update" -------------------------------------
# Plain `.Value = "May 2024"` gets auto-recognised as a date (month/year)
# and would change both the stored value (date serial) and the cell's
# style (a new numFmt gets attached). Round-tripping the text through a
# helper cell that is pre-formatted as Text, then pasting *values only*
# into the destination, keeps the destination cell's existing style/format
# untouched and stores a literal text value, just like the original cell.
$helper = $wsAbout.Range("Z1")
$helper.NumberFormat = "@"
$helper.Value = "May 2024"
$helper.Copy()
$wsAbout.Range("C6").PasteSpecial(-4163)   # xlPasteValues
$helper.Clear()

# --- REQUEST sheet (Table1) -------------------------------------------------
# Row for Seq 1 ("Account Type")
$wsRequest.Range("C3").Value = "Entity Value ( Contract No, Deal No, Reference ID, Account, etc)"

# Row for Seq 2 ("Account Number") + new Notes value
$wsRequest.Range("C4").Value = "Entity-Account Number"
$wsRequest.Range("F4").Value = "Model"

# Row for Seq 3 ("Account Check Digit") + new Notes value
$wsRequest.Range("C5").Value = "Entity-Account Check Digit"
$wsRequest.Range("F5").Value = "optional"

# --- RESPONSE sheet (Table2) -------------------------------------------------
# Row for Seq 2 ("Account Type")
$wsResponse.Range("C4").Value = "Entity Value ( Contract No, Deal No, Reference ID, Account, etc)"

# Row for Seq 3 ("Account Number")
$wsResponse.Range("C5").Value = "Entity-Account Number"

# Row for Seq 4 ("Account Check Digit")
$wsResponse.Range("C6").Value = "Entity-Account Check Digit"
